$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-8) are reshuffled: each target row receives the
# D, L, M, N, O, P, Q, S, T values that used to belong to a different
# source row (all other columns stay identical across these rows).
# Mapping: target row -> source row (values to copy in)
$rowData = @{
    2 = @{ D = 44210; L = "Primera"; M = 240; N = 15500; O = 16000; P = 15750; Q = "`$/caja 16 kilos granel"; S = 984; T = 16 }
    3 = @{ D = 44210; L = "Segunda"; M = 300; N = 12500; O = 13000; P = 12750; Q = "`$/caja 16 kilos granel"; S = 797; T = 16 }
    4 = @{ D = 44230; L = "Primera"; M = 160; N = 16500; O = 17000; P = 16750; Q = "`$/caja 18 kilos granel"; S = 931; T = 18 }
    5 = @{ D = 44230; L = "Segunda"; M = 160; N = 14500; O = 15000; P = 14750; Q = "`$/caja 18 kilos granel"; S = 819; T = 18 }
    6 = @{ D = 44224; L = "Especial"; M = 100; N = 16500; O = 17000; P = 16750; Q = "`$/caja 16 kilos granel"; S = 1047; T = 16 }
    7 = @{ D = 44224; L = "Primera"; M = 200; N = 14500; O = 15000; P = 14750; Q = "`$/caja 16 kilos granel"; S = 922; T = 16 }
    8 = @{ D = 44224; L = "Segunda"; M = 200; N = 12500; O = 13000; P = 12750; Q = "`$/caja 16 kilos granel"; S = 797; T = 16 }
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("Q$row").Value = $vals.Q
    $ws.Range("S$row").Value = $vals.S
    $ws.Range("T$row").Value = $vals.T
}
